# Add files via upload
# This script applies new log entries (row 19, date 42878) to the three
# worksheets "Ivan", "Fabian" and "Hannes", fills in a previously empty
# cell on "Hannes", and updates the selected cell per sheet.

$wb = $excel.ActiveWorkbook

$wsIvan   = $wb.Worksheets.Item("Ivan")
$wsFabian = $wb.Worksheets.Item("Fabian")
$wsHannes = $wb.Worksheets.Item("Hannes")

# --- Ivan (sheet1): new row 19 ---
$wsIvan.Range("A19").Value = 42878
$wsIvan.Range("A19").NumberFormat = "m/d/yyyy"
$wsIvan.Range("B19").Value = "push des gesamten Projektes, Fertigstellung der Dokumentation"
$wsIvan.Range("A19:B19").Borders.LineStyle = 1
$wsIvan.Range("B19").Interior.Pattern = -4142

# --- Fabian (sheet2): new row 19 ---
$wsFabian.Range("A19").Value = 42878
$wsFabian.Range("A19").NumberFormat = "m/d/yyyy"
$wsFabian.Range("B19").Value = "eventBearbeiten.php benutzerBearbeiten.php fertigstellen"
$wsFabian.Range("A19:B19").Borders.LineStyle = 1
$wsFabian.Range("B19").Interior.Pattern = -4142

# --- Hannes (sheet3): fill B18, add new row 19 ---
$wsHannes.Range("B18").Value = "Bug fix bei index.php"

$wsHannes.Range("A19").Value = 42878
$wsHannes.Range("A19").NumberFormat = "m/d/yyyy"
$wsHannes.Range("B19").Value = "Test auf Webhoster"
$wsHannes.Range("A19:B19").Borders.LineStyle = 1
$wsHannes.Range("B19").Interior.Pattern = -4142

# --- Update selections to match the authored file ---
$wsIvan.Range("B23").Select()
$wsFabian.Range("B18").Select()
$wsHannes.Range("B30").Select()
